# "SQL about the latest basic data update By Richard 20160112"
#
# The 商品 (Product) sheet's goodsId column (B3:B97) holds codes like
# "GD20160201000001" that were stored in the shared-string table with a
# stray trailing space (e.g. "GD20160201000001 "). This refresh of the
# basic data trims that trailing whitespace off every goods code in the
# column, which is exactly what a real Excel edit does: writing the
# trimmed text back creates brand-new shared-string entries (appended at
# the end of sharedStrings.xml) and re-points each cell at the new index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("商品")

$firstRow = 3
$lastRow = 97
$col = 2  # column B -> goodsId

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ($current -ne $null) {
        $trimmed = $current.Trim()
        if ($trimmed -ne $current) {
            $cell.Value = $trimmed
        }
    }
}
